# Weekly update: add a new week of price records (date serial 44522 = 2021-11-22)
# for "Terminal Hortofrutícola Agro Chillán - Limón". Two new rows are inserted
# right after the existing row 405, pushing all subsequent rows down by two.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("406:407").Insert()

$cols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T")

$row406 = @(
    7,
    "Terminal Hortofrutícola Agro Chillán",
    "Ñuble",
    44522,
    16,
    "Fruta",
    100102,
    "Cítricos",
    100102003,
    "Limón",
    "Sin especificar",
    "1a amarillo",
    200,
    6000,
    6500,
    6250,
    "`$/malla 16 kilos",
    "Región de O'Higgins",
    391,
    16
)

$row407 = @(
    7,
    "Terminal Hortofrutícola Agro Chillán",
    "Ñuble",
    44522,
    16,
    "Fruta",
    100102,
    "Cítricos",
    100102003,
    "Limón",
    "Sin especificar",
    "2a amarillo",
    120,
    5000,
    5500,
    5250,
    "`$/malla 16 kilos",
    "Región de O'Higgins",
    328,
    16
)

for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "406").Value() = $row406[$i]
    $ws.Range($cols[$i] + "407").Value() = $row407[$i]
}
